$d = $word.ActiveDocument

# 1. "welcome back to WICED WiFI 101." -> "... WiFi 101." (fix capitalization typo)
$d.Content.Find.Execute("WiFI", $true, $false, $false, $false, $false, $true, 1, $false, "WiFi", 2) | Out-Null

# 2. "Amazon Auroroa" -> "Amazon Aurora" (fix typo)
$d.Content.Find.Execute("Auroroa", $true, $false, $false, $false, $false, $true, 1, $false, "Aurora", 2) | Out-Null

# 3. "In fact one of" -> "In fact, one of" (add comma)
$d.Content.Find.Execute("In fact one of", $true, $false, $false, $false, $false, $true, 1, $false, "In fact, one of", 2) | Out-Null

# 4. ", aws will" -> ", AWS will" (capitalize AWS acronym)
$d.Content.Find.Execute(", aws will", $true, $false, $false, $false, $false, $true, 1, $false, ", AWS will", 2) | Out-Null

# 5. Move the "_GoBack" bookmark (Word's last-edit-location marker) to right after
#    the "AWS" we just typed, matching where the author's cursor ended up last.
#    Re-adding a bookmark with an existing name moves it from its old location.
$r = $d.Content
$r.Find.Execute(", AWS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
